# Update the "Förändrad" date column (C) for rows 2-26 from 45243 to 45244
# (i.e. advance the date by one day, from 2023-11-13 to 2023-11-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
